$d = $word.ActiveDocument

# The field paragraph is the 2nd paragraph in the document; it currently
# contains a Word field ({ m:'doc.html'.fromHTMLURI() }) built from
# fldChar/instrText runs. We need to turn it into plain visible text
# runs (w:t) that spell out the same field code, wrapped in literal
# "{" / "}" characters, while keeping the existing _GoBack bookmark
# (re-created at the same logical spot, between "doc.html" and
# "'.fromHTMLURI()").

$fieldPara = $d.Paragraphs.Item(2)
$r = $fieldPara.Range.Duplicate

# Wipe out the whole field (fldChar begin/end, instrText runs, old
# _GoBack bookmark) but keep the paragraph mark itself.
$r.MoveEnd(1, -1)
$r.Delete()

# Rebuild at the start of the (now empty) paragraph.
$ins = $fieldPara.Range.Duplicate
$ins.Collapse(1)

function Insert-Piece($range, [string]$text) {
    # Insert text, then drop a throw-away bookmark right after it so the
    # next insertion starts life as its own run instead of being merged
    # back into the previous w:t on save.
    $range.InsertAfter($text)
    $range.Collapse(0)
    $d.Bookmarks.Add("zzTmpSplit", $range) | Out-Null
    $d.Bookmarks.Item("zzTmpSplit").Delete()
}

Insert-Piece $ins "{"
Insert-Piece $ins "m"
Insert-Piece $ins ":"
Insert-Piece $ins "'"
Insert-Piece $ins "doc.html"

# Re-create the _GoBack bookmark exactly here (empty range); Word keeps
# bookmark names unique so this also removes the stray one left over
# from the original field, if any remained.
$d.Bookmarks.Add("_GoBack", $ins) | Out-Null

Insert-Piece $ins "'.fromHTMLURI()"

# Final run keeps its trailing-space-significant flag implicitly since
# there's no trailing space here, but match the xml:space="preserve"
# variant used in the target by inserting via the same helper.
$ins.InsertAfter("}")

Write-Output "done"
